# Shift the Timestamp column (A2:A97) forward by 9 days, from 2025-09-23
# to 2025-10-02, and zero out the two non-zero Actual Production values
# (B27, B28) that belonged to the old date, as part of adding GESS to the
# forecast portfolio.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 9
}

$ws.Range("B27").Value = 0
$ws.Range("B28").Value = 0
